# Update exercises data: fix "Back" -> "Upper back" for Barbell Row, and
# append four new exercises (rows 12-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 6 (Barbell Row): secondary/primary muscle label ---
$ws.Range("C6").Value = "Upper back"

# --- New row 12: Lateral raise ---
$ws.Range("A12").Value = "Lateral raise"
$ws.Range("B12").Value = "Dumbbell"
$ws.Range("C12").Value = "Side delts"
$ws.Range("D12").Value = "Traps"

# --- New row 13: Helm's row ---
$ws.Range("A13").Value = "Helm's row"
$ws.Range("B13").Value = "Dumbbell"
$ws.Range("C13").Value = "Upper back"
$ws.Range("D13").Value = "Lats"

# --- New row 14: Hammer curl ---
$ws.Range("A14").Value = "Hammer curl"
$ws.Range("B14").Value = "Bands"
$ws.Range("C14").Value = "Biceps"
$ws.Range("D14").Value = "Forearms"

# --- New row 15: Single arm pushdown (no variant/secondary muscle) ---
$ws.Range("A15").Value = "Single arm pushdown"
$ws.Range("C15").Value = "Triceps"

# --- Re-fit column widths for the new, wider content ---
$ws.Columns.Item(1).ColumnWidth = 17.498697916666668
$ws.Columns.Item(4).ColumnWidth = 15.830729166666666

# --- Update the active selection left by the editor ---
$ws.Range("C17").Select() | Out-Null
